# "separate dept from affiliations"
#
# - Sheet1 "PI hours" gains a new column F "app" holding the full list of
#   affiliations that used to live in column "dept", while "dept" itself is
#   narrowed down to each PI's single primary department.
# - Sheet2 "dept hours" is renamed to "department hours" and its data is
#   replaced with one row per primary department (ABE / ME / CS) using the
#   same hours/percentage numbers as each PI's row on sheet1.
# - A brand new sheet "unit(accumulative) hours" is appended, holding the
#   data that used to live on the old "dept hours" sheet (CSL/ABE/ME/AE/CS
#   accumulated hours), with its header relabeled "unit(accumulative)".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Keep references to the existing header / index-column formatting so it can
# be reused (via copy/paste of formats only) on newly written cells.
$headerStyleSource = $ws1.Range("B1")
$indexStyleSource = $ws1.Range("A2")

# ---------------------------------------------------------------------
# Sheet1 "PI hours": add the "app" column (old full affiliation lists)
# and shrink "dept" down to each PI's primary department.
# ---------------------------------------------------------------------
$headerStyleSource.Copy()
$ws1.Range("F1").PasteSpecial(-4122)
$ws1.Range("F1").Value = "app"

$ws1.Range("F2").Value = "['ABE', 'CSL']"
$ws1.Range("F3").Value = "['ME', 'AE', 'CSL']"
$ws1.Range("F4").Value = "['CS', 'CSL']"

$ws1.Range("E2").Value = "ABE"
$ws1.Range("E3").Value = "ME"
$ws1.Range("E4").Value = "CS"

# ---------------------------------------------------------------------
# Sheet2: rename "dept hours" -> "department hours" and replace its
# contents with the per-primary-department breakdown.
# ---------------------------------------------------------------------
$ws2.Name = "department hours"
$ws2.Cells.Clear()

$headerStyleSource.Copy()
$ws2.Range("B1:D1").PasteSpecial(-4122)
$ws2.Range("B1").Value = "dept"
$ws2.Range("C1").Value = "hours"
$ws2.Range("D1").Value = "percentage"

$indexStyleSource.Copy()
$ws2.Range("A2:A4").PasteSpecial(-4122)

$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "ABE"
$ws2.Range("C2").Value = 21.5
$ws2.Range("D2").Value = 51.19047619047619

$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "ME"
$ws2.Range("C3").Value = 14.5
$ws2.Range("D3").Value = 34.52380952380953

$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "CS"
$ws2.Range("C4").Value = 6
$ws2.Range("D4").Value = 14.28571428571429

# ---------------------------------------------------------------------
# New sheet3 "unit(accumulative) hours": carries over the data that used
# to be on the old "dept hours" sheet, with the header renamed.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "unit(accumulative) hours"

$headerStyleSource.Copy()
$ws3.Range("B1:D1").PasteSpecial(-4122)
$ws3.Range("B1").Value = "unit(accumulative)"
$ws3.Range("C1").Value = "hours"
$ws3.Range("D1").Value = "percentage"

$indexStyleSource.Copy()
$ws3.Range("A2:A6").PasteSpecial(-4122)

$ws3.Range("A2").Value = 0
$ws3.Range("B2").Value = "CSL"
$ws3.Range("C2").Value = 42
$ws3.Range("D2").Value = 42.63959390862944

$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = "ABE"
$ws3.Range("C3").Value = 21.5
$ws3.Range("D3").Value = 21.82741116751269

$ws3.Range("A4").Value = 2
$ws3.Range("B4").Value = "ME"
$ws3.Range("C4").Value = 14.5
$ws3.Range("D4").Value = 14.72081218274112

$ws3.Range("A5").Value = 3
$ws3.Range("B5").Value = "AE"
$ws3.Range("C5").Value = 14.5
$ws3.Range("D5").Value = 14.72081218274112

$ws3.Range("A6").Value = 4
$ws3.Range("B6").Value = "CS"
$ws3.Range("C6").Value = 6
$ws3.Range("D6").Value = 6.091370558375634

# Keep "PI hours" as the active/selected sheet, matching the original.
$ws1.Activate()
